# Regenerate merged AHB files
# 1. Rename header labels: "_old" -> "_FV2404", "_new" -> "_FV2410"
# 2. Turn the data range into an Excel Table ("Table1")
# 3. Freeze the header row (split/freeze panes at A2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $label = $cell.Value2
    if ($label -ne $null) {
        if ($label.EndsWith("_old")) {
            $cell.Value = $label.Substring(0, $label.Length - 4) + "_FV2404"
        } elseif ($label.EndsWith("_new")) {
            $cell.Value = $label.Substring(0, $label.Length - 4) + "_FV2410"
        }
    }
}

# Create an Excel Table ("ListObject") over the data range
$dataRange = $ws.Range("A1:U72")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# Freeze the header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
